# Update the "Game play description" worksheet (first sheet) with the new
# Nuclear Miami game-design text, replacing the old Rayman description.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the Wikipedia hyperlink that was attached to the old B2 cell; the
# new content for that cell is plain text with no link.
foreach ($hl in $ws.Hyperlinks) {
    $hl.Delete()
}

$ws.Range("B2").Value = "The game is a mix of Nuclear Throne with Hotline Miami, both of these are topdown fastpaced shooter games."
$ws.Range("B3").Value = "The player backstory is still no decided.`n"
$ws.Range("B4").Value = "In the game, for each level you spawn in an initial room with basic weaponry. He collects new weapons from the level as he progressed."
$ws.Range("B5").Value = "The exact ammount of health is not yet decided but it will be low, this so to provide a sense of challenge. "
$ws.Range("B6").Value = "At the end of each level there are stronger versions of the enemies, these are bosses. "
$ws.Range("B7").Value = "When it comes to power ups, it will mostly be health(med kits). Not decided yet."
$ws.Range("B8").Value = "There will be multiple different sorts of enemies, Melee, With simple weaponry, with more complex weaponry, and finally the boss with a rare weapon that you can take to the next level"
